$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Update the "进行中" (in progress) status on the previous week's block
#    (rows 86-89, column C) that was left blank before.
# ---------------------------------------------------------------------------
$ws.Range("C86").Value = "进行中"
$ws.Range("C87").Value = "进行中"
$ws.Range("C88").Value = "进行中"
$ws.Range("C89").Value = "进行中"

# ---------------------------------------------------------------------------
# 2. Write the values for the two new weekly-plan blocks appended at the
#    bottom of the sheet (rows 93-110).
# ---------------------------------------------------------------------------

# --- Block "2018.10.29" (rows 93-99) ---------------------------------------
$ws.Range("A93").Value = "日期：2018.10.29 第九周周一"

$ws.Range("A94").Value = "组员"
$ws.Range("B94").Value = "计划内容"
$ws.Range("C94").Value = "完成情况"
$ws.Range("D94").Value = "备注"

$ws.Range("A96").Value = "王嘉宇"
$ws.Range("B96").Value = "编写完善数据库"
$ws.Range("C96").Value = "进行中"

$ws.Range("A98").Value = "庞森杰"
$ws.Range("B98").Value = "使用安卓原生代码、百度地图sdk进行安卓端前端设计"
$ws.Range("C98").Value = "进行中"

$ws.Range("A100").Value = "总结："

# --- Block "2018.10.31" (rows 102-110) --------------------------------------
$ws.Range("A102").Value = "日期：2018.10.31 第九周周三"

$ws.Range("A103").Value = "组员"
$ws.Range("B103").Value = "计划内容"
$ws.Range("C103").Value = "完成情况"
$ws.Range("D103").Value = "备注"

$ws.Range("A104").Value = "余舒章"
$ws.Range("B104").Value = "编写完善数据库"
$ws.Range("C104").Value = "进行中"

$ws.Range("A105").Value = "王嘉宇"
$ws.Range("B105").Value = "编写完善数据库"
$ws.Range("C105").Value = "进行中"

$ws.Range("A106").Value = "许俊杰"
$ws.Range("B106").Value = "使用springmvc、mybatis进行框架搭建"
$ws.Range("C106").Value = "进行中"

$ws.Range("A107").Value = "庞森杰"
$ws.Range("B107").Value = "使用安卓原生代码、百度地图sdk进行安卓端前端设计"
$ws.Range("C107").Value = "进行中"

$ws.Range("A109").Value = "总结："

# ---------------------------------------------------------------------------
# 3. Copy the formatting from the matching existing blocks so the new rows
#    inherit the same styles (fonts/borders/alignment) used throughout the
#    rest of the sheet.
# ---------------------------------------------------------------------------
$ws.Range("A66:D72").Copy()
$ws.Range("A93:D99").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A82:D92").Copy()
$ws.Range("A100:D110").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4. Merge the title / summary cells for the new blocks.
# ---------------------------------------------------------------------------
$ws.Range("A93:D93").Merge()
$ws.Range("A100:D101").Merge()
$ws.Range("A102:D102").Merge()
$ws.Range("A109:D110").Merge()

# ---------------------------------------------------------------------------
# 5. Update the view so the sheet opens scrolled to the newly added rows.
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 73
$ws.Range("A109:D110").Select()
